# Apply edits to LOB1006.xlsx per target diff:
# 1. Insert a new row at 13 for the professor's name (moved from its
#    incorrect former homes at B10/C10 and B18/C18).
# 2. Correct several mis-placed cell values that resulted from the
#    original authoring mistakes (Objetivos, Programa resumido, Programa,
#    Metodo, Bibliografia).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new row for "Docentes responsaveis" value ---
$ws.Rows("13:13").Insert()

# The inserted row inherits a phantom styled A13 cell; clear it so column A
# stays empty on this row (matches target layout - only B13/C13 are used).
$ws.Range("A13").Clear()

# Fill in the professor's name on the new row.
$ws.Range("B13").Value = '6270264 - Juan Fernando Zapata Zapata'
$ws.Range("C13").Value = '6270264 - Juan Fernando Zapata Zapata'

# Copy the column B / C formatting down onto the new row 13 cells so they
# pick up the existing styles (wrap text, red font for C) instead of a
# freshly minted style.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Fix the "Objetivos:" row (10): was wrongly holding the professor's
#     name; should hold the PT objectives text. ---
$ws.Range("B10").Value = 'Familiarizar o aluno com os conceitos básicos de equações diferenciais e suas aplicações.'
$ws.Range("C10").Value = 'Familiarizar o aluno com os conceitos básicos de equações diferenciais e suas aplicações.'

# --- Fix "Programa resumido:" row (14, formerly 13): replace placeholder
#     "Semestral" with the actual short program summary. ---
$ws.Range("B14").Value = 'Sequencias e séries, equações diferenciais ordinárias de 1ª e 2ª ordem com aplicações, solução de equações diferenciais por series de potencia, Séries de Fourier e Problemas de valores de contorno.'
$ws.Range("C14").Value = 'Sequencias e séries, equações diferenciais ordinárias de 1ª e 2ª ordem com aplicações, solução de equações diferenciais por series de potencia, Séries de Fourier e Problemas de valores de contorno.'

# --- Fix "Programa:" row (16, formerly 15): replace the stray date value
#     with the full program description. ---
$ws.Range("B16").Value = 'Sequências e séries: Critérios de convergência, convergência condicional e absoluta, séries de potência, raio de convergência, derivação e integração termo a termo. Equações diferenciais ordinárias de 1ª e 2ª ordem: Equações exatas e não exatas, redução de ordem, Equação de Bernulli, método de variação de parâmetros e coeficientes a determinar, solução por séries de potencia de equações diferenciais, aplicações das equações diferenciais de 1ª e 2ª ordem.•Séries de Fourier: Teorema de convergência das séries de Fourier, Desigualdade de Bessel e Identidade de Parseval, equações em derivadas parciais e problemas de valores de contorno.'
$ws.Range("C16").Value = 'Sequências e séries: Critérios de convergência, convergência condicional e absoluta, séries de potência, raio de convergência, derivação e integração termo a termo. Equações diferenciais ordinárias de 1ª e 2ª ordem: Equações exatas e não exatas, redução de ordem, Equação de Bernulli, método de variação de parâmetros e coeficientes a determinar, solução por séries de potencia de equações diferenciais, aplicações das equações diferenciais de 1ª e 2ª ordem.•Séries de Fourier: Teorema de convergência das séries de Fourier, Desigualdade de Bessel e Identidade de Parseval, equações em derivadas parciais e problemas de valores de contorno.'

# --- Fix "Metodo:" row (19, formerly 18): was wrongly holding the
#     professor's name; should hold the evaluation method text. ---
$ws.Range("B19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("C19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'

# --- Fix "Bibliografia:" row (22, formerly 21): replace the recovery-norm
#     text (now correctly on row 21) with the actual bibliography text. ---
$ws.Range("B22").Value = '1.H. L. Guidorizzi, UM CURSO DE CÁLCULO, volume IV. Livros Técnicos e Científicos, 1987.2.BRANNAN, James R. BOYCE, W.E. Equações diferenciais: uma Introdução a métodos modernos e suas aplicações. Rio de Janeiro: LTC ED., 2008.3.ZILL, D.G. ; CULLEN, M.R. Equações Diferenciais São Paulo: Pearson Makron Books2006., v.1 e 2.4.W. Kaplan, CÁLCULO AVANÇADO, volume II, Edgard Blücher, São Paulo, 1972.5.BOYCE,W.E. ; DIPRIMA,R.C. Equações diferenciais e problemas de valores de contorno. 8.ed. Rio de Janeiro: LTC Editora, 2008.'
$ws.Range("C22").Value = '1.H. L. Guidorizzi, UM CURSO DE CÁLCULO, volume IV. Livros Técnicos e Científicos, 1987.2.BRANNAN, James R. BOYCE, W.E. Equações diferenciais: uma Introdução a métodos modernos e suas aplicações. Rio de Janeiro: LTC ED., 2008.3.ZILL, D.G. ; CULLEN, M.R. Equações Diferenciais São Paulo: Pearson Makron Books2006., v.1 e 2.4.W. Kaplan, CÁLCULO AVANÇADO, volume II, Edgard Blücher, São Paulo, 1972.5.BOYCE,W.E. ; DIPRIMA,R.C. Equações diferenciais e problemas de valores de contorno. 8.ed. Rio de Janeiro: LTC Editora, 2008.'
